$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos price/volume table refresh (GitHub Actions data pull).
# Write each cell with a leading apostrophe (Excel text-prefix) so
# numeric-looking strings like "1.00" / "0.000270" stay text and keep
# their exact formatting instead of being coerced to numbers, then
# reset the style back to Normal so no quotePrefix style sticks and
# the cell style index is left untouched.
function Set-TextCell($addr, $val) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextCell "D2" '71.048.59'
Set-TextCell "E2" '  -1.56%  '
Set-TextCell "D3" '3.952.90'
Set-TextCell "E3" '  -2.09%  '
Set-TextCell "D4" '1.00'
Set-TextCell "E4" '  -0.08%  '
Set-TextCell "D5" '537.78'
Set-TextCell "E5" '  +3.18%  '
Set-TextCell "D6" '148.12'
Set-TextCell "E6" '  +0.65%  '
Set-TextCell "D7" '3.948.11'
Set-TextCell "E7" '  -1.95%  '
Set-TextCell "E8" '  -5.95%  '
Set-TextCell "D9" '1.00'
Set-TextCell "E9" '  -0.03%  '
Set-TextCell "D10" '0.737'
Set-TextCell "E10" '  -4.92%  '
Set-TextCell "E11" '  -5.82%  '
Set-TextCell "D12" '55.26'
Set-TextCell "E12" '  +13.55%  '
Set-TextCell "E13" '  -3.66%  '
Set-TextCell "D14" '10.63'
Set-TextCell "E14" '  -4.66%  '
Set-TextCell "D15" '4.579.91'
Set-TextCell "E15" '  -2.10%  '
Set-TextCell "D16" '3.947.26'
Set-TextCell "E16" '  -2.57%  '
Set-TextCell "D17" '20.59'
Set-TextCell "D18" '13.81'
Set-TextCell "E18" '  -2.95%  '
Set-TextCell "E19" '  -1.44%  '
Set-TextCell "E20" '  -4.13%  '
Set-TextCell "D21" '70.937.39'
Set-TextCell "E21" '  -1.67%  '
Set-TextCell "D22" '425.44'
Set-TextCell "E22" '  -4.19%  '
Set-TextCell "E23" '  +0.11%  '
Set-TextCell "D24" '96.98'
Set-TextCell "E24" '  -7.39%  '
Set-TextCell "D25" '4.22'
Set-TextCell "E25" '  +4.81%  '
Set-TextCell "D26" '14.41'
Set-TextCell "E26" '  -4.09%  '
Set-TextCell "D27" '11.38'
Set-TextCell "E27" '  -1.25%  '
Set-TextCell "D28" '10.63'
Set-TextCell "E28" '  -3.97%  '
Set-TextCell "D29" '3.76'
Set-TextCell "D30" '5.91'
Set-TextCell "E30" '  +1.46%  '
Set-TextCell "D31" '36.38'
Set-TextCell "E31" '  -4.01%  '
Set-TextCell "D32" '7.81'
Set-TextCell "E32" '  +15.64%  '
Set-TextCell "D33" '50.08'
Set-TextCell "E33" '  +17.72%  '
Set-TextCell "D34" '0.130'
Set-TextCell "E34" '  -0.20%  '
Set-TextCell "E35" '  -2.99%  '
Set-TextCell "D36" '683.98'
Set-TextCell "E36" '  +1.23%  '
Set-TextCell "D37" '65.06'
Set-TextCell "E37" '  -3.10%  '
Set-TextCell "E38" '  +2.79%  '
Set-TextCell "D39" '0.0₃0819'
Set-TextCell "E39" '  -5.50%  '
Set-TextCell "E40" '  -1.77%  '
Set-TextCell "D41" '3.39'
Set-TextCell "E41" '  -3.57%  '
Set-TextCell "E42" '  +0.06%  '
Set-TextCell "E43" '  +0.02%  '
Set-TextCell "D44" '0.0481'
Set-TextCell "E44" '  -4.15%  '
Set-TextCell "E45" '  -2.25%  '
Set-TextCell "E46" '  -7.58%  '
Set-TextCell "B47" 'THORChain'
Set-TextCell "C47" 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextCell "D47" '9.90'
Set-TextCell "E47" '  +6.45%  '
Set-TextCell "B48" 'Fetch.AI'
Set-TextCell "C48" 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextCell "D48" '2.69'
Set-TextCell "E48" '  -0.90%  '
Set-TextCell "E49" '  -6.57%  '
Set-TextCell "D50" '3.00'
Set-TextCell "E50" '  -2.05%  '
Set-TextCell "D51" '0.000270'
Set-TextCell "E51" '  -1.72%  '
